# "some tweaks, levels 3-5"
# Updates a couple of existing shape-description strings, inserts a new
# "tri_desc" key/value row right after "poly_desc", and appends six new
# localization rows (level_intro_3_0 .. level_intro_5_0) for the new
# triangle levels (3-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing descriptions (ellipse_desc / poly_desc) -------------
$ws.Range("B68").Value = '· Round, has a center point.\n\n· Dimensions defined by two axis.'
$ws.Range("B70").Value = '· Formed by three or more straight lines (sides) connected in a loop.\n\n· Points are plotted on a plane.'

# --- Insert a new row for "tri_desc" right after "poly_desc" (row 70) ----
$ws.Rows("71:71").Insert()
$ws.Range("A71").Value = "tri_desc"
$ws.Range("B71").Value = '· Formed by three straight lines connected in a loop.\n· Points are plotted on a plane.\n· Angles sum up to 180°.'

# --- Append new level-intro rows for levels 3-5 at the bottom ------------
$ws.Range("A85").Value = "level_intro_3_0"
$ws.Range("B85").Value = 'Now that we''ve covered some of the polygons, let''s take a closer look at triangles.'

$ws.Range("A86").Value = "level_intro_3_1"
$ws.Range("B86").Value = 'As a sub-category of polygons, triangles follow the same attribute of being formed by a number of straight lines connected in a loop. In this case, there are three.'

$ws.Range("A87").Value = "level_intro_3_2"
$ws.Range("B87").Value = 'With that in mind, we can then further identify sub-categories under the triangle with additional attributes.'

$ws.Range("A88").Value = "level_intro_3_3"
$ws.Range("B88").Value = 'For this level, we will be looking at different triangles based on their angle values. Be sure to remember them!'

$ws.Range("A89").Value = "level_intro_4_0"
$ws.Range("B89").Value = 'Now that we''ve seen some triangles based on their angle values, let''s take a look at the ones based on their side values.'

$ws.Range("A90").Value = "level_intro_5_0"
$ws.Range("B90").Value = 'For this level, we will be categorizing triangles based on their angles, and side lengths. These triangles will have more than one category that fit their attributes.'

# --- Match the view state captured in the saved workbook -----------------
$null = $ws.Range("A73").Select()
$excel.ActiveWindow.ScrollRow = 73
$null = $ws.Range("B90").Select()
